# Prefix each "command" name (column A) in the step/variant sheets with
# the protocol (sheet) name, per commit: "fix: unique command names in XLSX -
# prefix protocol name to each step".
#
# Sheets MaxJourney, NRWaves, PersonalMax, PositiveSpin, ReEngagement are
# left untouched; every other sheet (price1, price2, discount1, discount2,
# free1, free2, nomoney1, nomoney2, noppv1, noppv2, card1, card2, nosex1,
# nosex2, offtopic1, offtopic2, real1, real2, voice1, voice2, customyes1,
# customyes2, customno1, customno2, done1, done2, cumcontrol, dickpic,
# boosters) gets its A-column values (rows below the header) prefixed with
# "<SheetName> ".

$wb = $excel.ActiveWorkbook

$skipSheets = @("MaxJourney", "NRWaves", "PersonalMax", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name

    if ($skipSheets -contains $name) {
        continue
    }

    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2

        if ($current -eq $null) {
            continue
        }

        $currentText = [string]$current

        if ($currentText -eq "") {
            continue
        }

        $prefix = $name + " "
        if ($currentText.StartsWith($prefix)) {
            continue
        }

        $cell.Value = $prefix + $currentText
    }
}
